$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1039
$ws1.Range("F11").Value = 492
$ws1.Range("F13").Value = 151
$ws1.Range("F14").Value = 12344
$ws1.Range("F15").Value = 89
$ws1.Range("F16").Value = 5473

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 116

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 116
$ws4.Range("F7").Value = 1039
$ws4.Range("F13").Value = 492
$ws4.Range("F15").Value = 151
$ws4.Range("F16").Value = 12344
$ws4.Range("F18").Value = 89
$ws4.Range("F19").Value = 5473
